$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Meter -> Metre, fix LaTeX escaping (double backslash -> single backslash)
$ws.Range("A2").Value = "c_\mathit{m}"
$ws.Range("C2").Value = "Metre"

# Row 3: Second stays, fix LaTeX escaping
$ws.Range("A3").Value = "c_\mathit{s}"

# Row 4: now Joule (was Lspeed); clear old formula in B4
$ws.Range("A4").Value = "c_\mathit{J}"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "Joule"

# Row 5: now Kelvin (was Planck length); clear old formula in B5
$ws.Range("A5").Value = "c_\mathit[K}"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "Kelvin"

# Row 6: now Gramm (was Planck constant); clear old formula in B6
$ws.Range("A6").Value = "c_\mathit{g]"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "Gramm"

# Row 7: new - Newton
$ws.Range("A7").Value = "c_\mathit{N]"
$ws.Range("C7").Value = "Newton"

# Row 8: new - Kilogramm
$ws.Range("A8").Value = "c_\mathit{kg]"
$ws.Range("B8").Value = "\mathit{g\cdot1000}"
$ws.Range("C8").Value = "Kilogramm"

# Rows 9-11 intentionally left blank

# Row 12: Lspeed (moved down from row 4, with corrected LaTeX)
$ws.Range("A12").Value = "c_\mathit{c}"
$ws.Range("B12").Value = "(299792458 \cdot \frac{c_\mathit{m}}{c_\mathit{s}})"
$ws.Range("C12").Value = "Lspeed"

# Row 13: Planck length (moved down from row 5, with corrected LaTeX)
$ws.Range("A13").Value = "c_\mathit{L}"
$ws.Range("B13").Value = "(1.616255 \cdot 10^{-35} \cdot c_\mathit{m})"
$ws.Range("C13").Value = "Planck length"

# Row 14: Planck constant (moved down from row 6, with corrected LaTeX)
$ws.Range("A14").Value = "c_\mathit{h}"
$ws.Range("B14").Value = "(1.055 \cdot 10^{-34} \cdot c_\mathit{J} \dot c_\mathit{s})"
$ws.Range("C14").Value = "Planck constant"

# Row 15: new - Reduced Planck constant
$ws.Range("A15").Value = "c_\mathit{rh}"
$ws.Range("B15").Value = "\frac{\mathit{h}}{\mathit{2\cdot pi}}"
$ws.Range("C15").Value = "Reduced Planck constant"

# Row 16: new - Gravitational constant
$ws.Range("A16").Value = "c_\mathit{G}"
$ws.Range("B16").Value = "6.674\cdot10^{-11}\cdot\frac{N\cdotm^2}{kg^2}"
$ws.Range("C16").Value = "Gravitational constant"

$ws.Range("C16").Select()
